## Applies the commit's changes:
##  1. Adds a new worksheet "ODI Batting Extra" (after "ODI Batting") with
##     per-match aggregate stats.
##  2. Clears the stray empty cell B12 on the "ODI Batting" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Clear B12 on "ODI Batting" -> the cell becomes genuinely empty (no
#    cell record at all), matching the removal of the empty <c r="B12"/>.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B12").ClearContents()
$odiBatting.Range("E12").Value2 = " "

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Batting Extra" worksheet right after "ODI Batting".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Match the look & feel of the sibling sheets (outline summary flags +
# page margins used throughout the rest of the workbook).
$newSheet.Outline.SummaryBelow = $true
$newSheet.Outline.SummaryRight = $true
$newSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.TopMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Header row - grab the bold/centered/bordered look of the other sheets'
# header row by copying its formatting over, then stamp in the real text.
$odiBatting.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# MATCH_CODE (A), NUM_4 (C), NUM_6 (D) and PERCENT_RUNS_OF_TOTAL (E) are
# stored as plain text in the source data (e.g. "4452", "2", "20.45%") -
# force a text format on those columns so Excel doesn't silently convert
# them to numbers/percentages.
$textCols = @(1, 3, 4, 5)
foreach ($col in $textCols) {
    $colRange = $newSheet.Range($newSheet.Cells.Item(2, $col), $newSheet.Cells.Item(19, $col))
    $colRange.NumberFormat = "@"
}

# Per-match data: MatchCode, BattingPosition, Num4, Num6, PercentRuns, ManOfMatch
$data = @(
    @("4452", 3,    "2",  "0", "20.45%", "NO"),
    @("4453", 3,    "7",  "0", "26.18%", "NO"),
    @("4455", 3,    "17", "0", "39.62%", "YES"),
    @("4636", 4,    "0",  "0", "2.11%",  "NO"),
    @("4639", 3,    "0",  "0", "2.83%",  "NO"),
    @("4642", $null, $null, $null, $null, "NO"),
    @("4647", $null, $null, $null, $null, "NO"),
    @("4648", 2,    "0",  "0", "6.10%",  "NO"),
    @("4649", 2,    "5",  "0", "8.68%",  "NO"),
    @("4669", 2,    "3",  "0", "7.77%",  "NO"),
    @("4673", $null, $null, $null, $null, "NO"),
    @("4676", 2,    "6",  "0", "36.54%", "NO"),
    @("4686", $null, $null, $null, $null, "NO"),
    @("4688", 2,    "13", "1", "38.70%", "YES"),
    @("4690", $null, $null, $null, $null, "NO"),
    @("4692", $null, $null, $null, $null, "NO"),
    @("4695", 2,    "1",  "0", "6.48%",  "NO"),
    @("4697", 2,    "12", "8", "46.78%", "NO")
)

$row = 2
foreach ($rec in $data) {
    $newSheet.Cells.Item($row, 1).Value2 = $rec[0]
    if ($null -ne $rec[1]) { $newSheet.Cells.Item($row, 2).Value2 = $rec[1] }
    if ($null -ne $rec[2]) { $newSheet.Cells.Item($row, 3).Value2 = $rec[2] }
    if ($null -ne $rec[3]) { $newSheet.Cells.Item($row, 4).Value2 = $rec[3] }
    if ($null -ne $rec[4]) { $newSheet.Cells.Item($row, 5).Value2 = $rec[4] }
    $newSheet.Cells.Item($row, 6).Value2 = $rec[5]
    $row++
}

# Leave the workbook's active sheet untouched (still "Player Info"),
# matching the original activeTab="0".
$wb.Worksheets.Item("Player Info").Activate()
